$d = $word.ActiveDocument
$CR = [char]13

# ------------------------------------------------------------------
# Part 1: insert a new "Meta description" paragraph right after the
# Heading1 title paragraph.
# ------------------------------------------------------------------

$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

# Locate the existing bold "Play Crystal Crush for Free - Online Slot
# Game Review" paragraph further down the document (its run structure,
# including the leading empty run, is exactly what we want to reuse)
# and copy its formatted text into the freshly inserted paragraph.
# (Paragraph.Range.Text always carries a trailing paragraph mark, so
# the comparison target needs one too.)
$boldTarget = "Play Crystal Crush for Free - Online Slot Game Review" + $CR
$boldSourcePara = $null
$count = $d.Paragraphs.Count
for ($i = 2; $i -le $count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -eq $boldTarget) {
        $boldSourcePara = $cand
        break
    }
}
if ($boldSourcePara -eq $null) {
    throw "Could not locate the source bold title paragraph to copy formatting from."
}
$metaPara.Range.FormattedText = $boldSourcePara.Range.FormattedText

# Turn the copied bold text into "Meta description".
$boldStart = $metaPara.Range.Start
$boldEnd = $boldStart + "Play Crystal Crush for Free - Online Slot Game Review".Length
$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Text = "Meta description"

# Append the (non-bold) remainder of the sentence at the end of the
# paragraph, just before its paragraph mark.
$pEnd = $metaPara.Range.End
$tail = $d.Range($pEnd - 1, $pEnd - 1)
$tail.InsertAfter(": Play Crystal Crush for free and read our review. Crystal Crush features original gameplay and a hexagon-shaped grid that makes it stand out.")

# ------------------------------------------------------------------
# Part 2: near the end of the document, drop the duplicate bold
# "Play Crystal Crush for Free - Online Slot Game Review" paragraph
# and replace the italic paragraph's text with the new image prompt.
# ------------------------------------------------------------------

$boldDupTarget = "Play Crystal Crush for Free - Online Slot Game Review" + $CR
$italicTarget = "Play Crystal Crush for free and read our review. Crystal Crush features original gameplay and a hexagon-shaped grid that makes it stand out." + $CR

$count = $d.Paragraphs.Count
$boldDupIndex = -1
for ($i = 3; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -eq $boldDupTarget) {
        $boldDupIndex = $i
        break
    }
}

if ($boldDupIndex -eq -1) {
    throw "Could not locate the duplicate bold title paragraph to delete."
}
$d.Paragraphs($boldDupIndex).Range.Delete()

# Recompute the italic paragraph's index after the deletion shifted
# everything down by one.
$count = $d.Paragraphs.Count
$italicIndex = -1
for ($i = 3; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -eq $italicTarget) {
        $italicIndex = $i
        break
    }
}

if ($italicIndex -eq -1) {
    throw "Could not locate the italic paragraph whose text needs replacing."
}
$italicPara = $d.Paragraphs($italicIndex)
$start = $italicPara.Range.Start
$end = $italicPara.Range.End
$textRange = $d.Range($start, $end - 1)
$textRange.Text = "Prompt: Design a feature image for Crystal Crush that portrays a happy Maya warrior with glasses in cartoon style. DALLE, please design a feature image for Crystal Crush that captures the essence of this innovative and exciting slot game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be depicted in a tropical setting with crystals and gems surrounding him. The colors should be bright and vibrant, capturing the game's fun and playful nature. The image should be eye-catching and convey the unique mechanics of the game that set it apart from other online slots. Let your creativity run wild and provide an image that will make players want to dive right into the world of Crystal Crush."

Write-Output "done"
